$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions): update "想去人数" (want-to-go count) values
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 1214
$wsExhibit.Range("F4").Value = 2688

# Sheet "全部类型" (All types): same two events appear here, update matching cells
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 1214
$wsAll.Range("F6").Value = 2688
